$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$shp = $master.Shapes.Item(3)
$tr2 = $shp.TextFrame2.TextRange
Write-Host "Chars Count=" $tr2.Characters.Count
$ch = $tr2.Characters(1, $tr2.Length)
Write-Host "ch Text=[" $ch.Text "]"
$ch.Text = "2024. 04. 16."
Write-Host "Text after=[" $tr2.Text "]"
